{"js": "const replacements = [\n  [\"370\u00f78=46, 2\", \"850\u00f72=425, 0\"],\n  [\"464\u00f77=66, 2\", \"984\u00f79=109, 3\"],\n  [\"489\u00f73=163, 0\", \"392\u00f75=78, 2\"],\n  [\"914\u00f76=152, 2\", \"931\u00f78=116, 3\"],\n  [\"865\u00f77=123, 4\", \"648\u00f76=108, 0\"],\n  [\"274\u00f74=68, 2\", \"792\u00f79=88, 0\"],\n  [\"685\u00f74=171, 1\", \"915\u00f77=130, 5\"],\n  [\"384\u00f76=64, 0\", \"271\u00f74=67, 3\"],\n  [\"207\u00f72=103, 1\", \"161\u00f79=17, 8\"],\n  [\"816\u00f77=116, 4\", \"869\u00f76=144, 5\"],\n  [\"605\u00f75=121, 0\", \"788\u00f75=157, 3\"],\n  [\"853\u00f76=142, 1\", \"825\u00f73=275, 0\"],\n  [\"217\u00f79=24, 1\", \"362\u00f73=120, 2\"],\n  [\"662\u00f78=82, 6\", \"675\u00f74=168, 3\"],\n  [\"828\u00f77=118, 2\", \"573\u00f74=143, 1\"],\n  [\"483\u00f74=120, 3\", \"310\u00f74=77, 2\"],\n  [\"107\u00f73=35, 2\", \"660\u00f79=73, 3\"],\n  [\"455\u00f75=91, 0\", \"125\u00f77=17, 6\"],\n  [\"968\u00f72=484, 0\", \"623\u00f76=103, 5\"],\n  [\"528\u00f78=66, 0\", \"964\u00f73=321, 1\"],\n  [\"394\u00f78=49, 2\", \"265\u00f74=66, 1\"],\n  [\"230\u00f77=32, 6\", \"433\u00f75=86, 3\"],\n  [\"804\u00f72=402, 0\", \"472\u00f75=94, 2\"],\n  [\"154\u00f77=22, 0\", \"527\u00f77=75, 2\"],\n  [\"402\u00f78=50, 2\", \"793\u00f73=264, 1\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"370\u00f78=46, 2\", \"850\u00f72=425, 0\"),\n    @(\"464\u00f77=66, 2\", \"984\u00f79=109, 3\"),\n    @(\"489\u00f73=163, 0\", \"392\u00f75=78, 2\"),\n    @(\"914\u00f76=152, 2\", \"931\u00f78=116, 3\"),\n    @(\"865\u00f77=123, 4\", \"648\u00f76=108, 0\"),\n    @(\"274\u00f74=68, 2\", \"792\u00f79=88, 0\"),\n    @(\"685\u00f74=171, 1\", \"915\u00f77=130, 5\"),\n    @(\"384\u00f76=64, 0\", \"271\u00f74=67, 3\"),\n    @(\"207\u00f72=103, 1\", \"161\u00f79=17, 8\"),\n    @(\"816\u00f77=116, 4\", \"869\u00f76=144, 5\"),\n    @(\"605\u00f75=121, 0\", \"788\u00f75=157, 3\"),\n    @(\"853\u00f76=142, 1\", \"825\u00f73=275, 0\"),\n    @(\"217\u00f79=24, 1\", \"362\u00f73=120, 2\"),\n    @(\"662\u00f78=82, 6\", \"675\u00f74=168, 3\"),\n    @(\"828\u00f77=118, 2\", \"573\u00f74=143, 1\"),\n    @(\"483\u00f74=120, 3\", \"310\u00f74=77, 2\"),\n    @(\"107\u00f73=35, 2\", \"660\u00f79=73, 3\"),\n    @(\"455\u00f75=91, 0\", \"125\u00f77=17, 6\"),\n    @(\"968\u00f72=484, 0\", \"623\u00f76=103, 5\"),\n    @(\"528\u00f78=66, 0\", \"964\u00f73=321, 1\"),\n    @(\"394\u00f78=49, 2\", \"265\u00f74=66, 1\"),\n    @(\"230\u00f77=32, 6\", \"433\u00f75=86, 3\"),\n    @(\"804\u00f72=402, 0\", \"472\u00f75=94, 2\"),\n    @(\"154\u00f77=22, 0\", \"527\u00f77=75, 2\"),\n    @(\"402\u00f78=50, 2\", \"793\u00f73=264, 1\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
